$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($row, $startCol, $values) {
    $n = $values.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $startColIdx = [int][char]$startCol - [int][char]'A' + 1
    $endColIdx = $startColIdx + $n - 1
    $startAddr = "$startCol$row"
    $endCol = [char]([int][char]'A' + $endColIdx - 1)
    $endAddr = "$endCol$row"
    $rangeAddr = "$startAddr`:$endAddr"
    $ws.Range($rangeAddr).Value = $arr
}

Set-RowValues 2 "B" @(1.02, 1.02695203671007, 1.028328299136783, 1.035237822065761, 1.041803740305805)
Set-RowValues 2 "I" @(1.025334150679314, 1.032112615282299, 1.031145435481519, 1.038035001126089, 1.044582220857773, 1.033578333580972)
Set-RowValues 3 "B" @(1.02, 1.028425232137423, 1.029631961451926, 1.036650300617254, 1.043557932578313)
Set-RowValues 3 "I" @(1.025536983666918, 1.033223097584837, 1.032255673016592, 1.039255229838042, 1.046144642233662, 1.034690392895753)
Set-RowValues 4 "B" @(1.02, 1.029370843346969, 1.03046877032298, 1.037554053670829, 1.044673777372807)
Set-RowValues 4 "I" @(1.025657876763042, 1.033933687660707, 1.032966803884295, 1.040034094668599, 1.047136003077589, 1.035401992091036)
Set-RowValues 5 "B" @(1.02, 1.029766572281228, 1.030818970224196, 1.037931572896479, 1.045138317271887)
Set-RowValues 5 "I" @(1.025706230056725, 1.034230533349913, 1.033264043358741, 1.040358992102016, 1.047548112421339, 1.03569925933517)
Set-RowValues 6 "B" @(1.02, 1.029832911684563, 1.030877677410525, 1.037994818964872, 1.045216049308911)
Set-RowValues 6 "I" @(1.025714204161799, 1.034280265010944, 1.033313850871799, 1.040413395729891, 1.047617035337502, 1.035749061620864)
Set-RowValues 7 "B" @(1.02, 1.029376138160909, 1.030473455949138, 1.037559107569086, 1.044680002443213)
Set-RowValues 7 "I" @(1.025658532557239, 1.033937661512218, 1.032970782344829, 1.040038445900702, 1.047141527955253, 1.035405971585872)
Set-RowValues 8 "B" @(1.02, 1.027451507823952, 1.028770288034655, 1.035717307604062, 1.042400587639863)
Set-RowValues 8 "I" @(1.025404848216088, 1.032489570303923, 1.031522162420935, 1.038449615395558, 1.045114338450558, 1.033955823921963)
Set-RowValues 9 "B" @(1.02, 1.0240003671133, 1.02571641300826, 1.032392296412668, 1.038234551398291)
Set-RowValues 9 "I" @(1.02487809101332, 1.02987584451053, 1.028912906184225, 1.035566688072133, 1.041389856971178, 1.031338386338145)
Set-RowValues 10 "B" @(1.02, 1.021657774913725, 1.023643603798699, 1.030120233343479, 1.03535359094919)
Set-RowValues 10 "I" @(1.02447261189647, 1.028090195721257, 1.027133939128689, 1.033586972012599, 1.038801500657834, 1.029550201722775)
Set-RowValues 11 "B" @(1.02, 1.020633080685022, 1.022736950987595, 1.029122805664956, 1.03408077023702)
Set-RowValues 11 "I" @(1.024283973241851, 1.027306397382977, 1.026353929464155, 1.032715592772162, 1.037655003309166, 1.028765290300951)
Set-RowValues 12 "B" @(1.02, 1.020250874109158, 1.022398778611398, 1.028750229881936, 1.033604110529106)
Set-RowValues 12 "I" @(1.024211924498206, 1.027013634870503, 1.026062710368655, 1.032389757847781, 1.037225212800216, 1.028472112032129)
Set-RowValues 13 "B" @(1.02, 1.02033293125313, 1.022471381671761, 1.028830243781061, 1.033706532273989)
Set-RowValues 13 "I" @(1.024227469088018, 1.027076507490683, 1.026125245733401, 1.032459749223683, 1.037317583279037, 1.028535073938642)
Set-RowValues 14 "B" @(1.02, 1.020601520013962, 1.022709026327905, 1.029092051254627, 1.034041448911041)
Set-RowValues 14 "I" @(1.024278058165635, 1.02728223088289, 1.026329887756639, 1.032688703578865, 1.037619557264316, 1.028741089481663)
Set-RowValues 15 "B" @(1.02, 1.020766794663387, 1.022855260391097, 1.029253081671567, 1.034247286164951)
Set-RowValues 15 "I" @(1.024308964862676, 1.027408767544069, 1.026455776239605, 1.032829481676995, 1.037805090559311, 1.028867805839414)
Set-RowValues 16 "B" @(1.02, 1.021725559636049, 1.023703580663686, 1.030186139057569, 1.035437524014317)
Set-RowValues 16 "I" @(1.024484854543758, 1.028141987774226, 1.027185498740502, 1.033644501100336, 1.038877042436684, 1.029602067326401)
Set-RowValues 17 "B" @(1.02, 1.022324174951516, 1.024233247676435, 1.030767748665445, 1.036177295095913)
Set-RowValues 17 "I" @(1.024591676234832, 1.028599057088119, 1.027640614857939, 1.034151925468072, 1.039542517653138, 1.030059785731145)
Set-RowValues 18 "B" @(1.02, 1.022672342396025, 1.024541316488006, 1.031105681562994, 1.036606350447978)
Set-RowValues 18 "I" @(1.024652724325991, 1.028864637349602, 1.027905142135981, 1.034446535158773, 1.039928200943681, 1.030325743147061)
Set-RowValues 19 "B" @(1.02, 1.022790890853636, 1.024646212184933, 1.031220686965702, 1.036752235313707)
Set-RowValues 19 "I" @(1.024673327068774, 1.028955021207597, 1.027995181541826, 1.03454675955322, 1.040059290751804, 1.030416255360503)
Set-RowValues 20 "B" @(1.02, 1.022260052343485, 1.024176510357532, 1.030705483271037, 1.036098177593963)
Set-RowValues 20 "I" @(1.024580345643575, 1.028550123663682, 1.027591882057043, 1.03409762483436, 1.039471375148229, 1.030010782815633)
Set-RowValues 21 "B" @(1.02, 1.02052247149912, 1.022639084843775, 1.029015013431654, 1.0339429319215)
Set-RowValues 21 "I" @(1.024263215745136, 1.027221695594264, 1.026269667132859, 1.032621342349305, 1.03753074249632, 1.028680468225983)
Set-RowValues 22 "B" @(1.02, 1.019420769627559, 1.021664320733623, 1.027940052341917, 1.032565372695479)
Set-RowValues 22 "I" @(1.024052358661788, 1.026377041723628, 1.025429708213921, 1.031680592447062, 1.036287814640317, 1.027834614849968)
Set-RowValues 23 "B" @(1.02, 1.02000568886333, 1.022181842781704, 1.028511070696013, 1.033297797063313)
Set-RowValues 23 "I" @(1.024165230868542, 1.026825712689849, 1.025875814994903, 1.032180505558566, 1.036948895999379, 1.028283922980424)
Set-RowValues 24 "B" @(1.02, 1.022289029662023, 1.02420215021819, 1.030733622368928, 1.03613393492542)
Set-RowValues 24 "I" @(1.02458546934014, 1.028572237723803, 1.027613905197503, 1.034122165161121, 1.03950352903489, 1.030032928280256)
Set-RowValues 25 "B" @(1.02, 1.024899812923783, 1.026512302106446, 1.033261506263614, 1.039329589526382)
Set-RowValues 25 "I" @(1.025023786491368, 1.030559049762864, 1.029594308941602, 1.032022561820302, 1.042371067472507, 1.032022561820302)
